$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Actual" burndown values in column C (rows 4-9).
$ws.Range("C4").Value = 19
$ws.Range("C5").Value = 18.5
$ws.Range("C6").Value = 18
$ws.Range("C7").Value = 17.5
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 4

# Update the active selection to match the saved view (D9).
$ws.Range("D9").Select() | Out-Null
